$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A1: "judul_proposal" -> "judul"
$ws.Range("A1").Value = "judul"
